$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 45000
$ws.Range("J3").Value = 45000
$ws.Range("L3").Value = 45000
$ws.Range("N3").Value = -45228

$ws.Range("H39").Value = 184
$ws.Range("I39").Value = 141.05263
$ws.Range("J39").Value = 1000
$ws.Range("K39").Value = 423.15789
$ws.Range("L39").Value = 3000
$ws.Range("M39").Value = -127.15789
$ws.Range("N39").Value = -3592

$ws.Range("H40").Value = 1997.5
$ws.Range("J40").Value = 1997.5
$ws.Range("L40").Value = 1997.5
$ws.Range("N40").Value = -2347.5

$ws.Range("H53").Value = 356.9091
$ws.Range("I53").Value = 340.625
$ws.Range("K53").Value = 340.625
$ws.Range("M53").Value = 296.375

$ws.Range("H55").Value = 394
$ws.Range("I55").Value = 308.33334
$ws.Range("K55").Value = 308.33334
$ws.Range("M55").Value = -94.33334000000002

$ws.Range("H62").Value = 5091.857
$ws.Range("I62").Value = 4488.5
$ws.Range("J62").Value = 5896.3335
$ws.Range("K62").Value = 4488.5
$ws.Range("L62").Value = 5896.3335
$ws.Range("M62").Value = -3864.5
$ws.Range("N62").Value = -7144.3335

$ws.Range("H65").Value = 5091.857
$ws.Range("I65").Value = 4488.5
$ws.Range("J65").Value = 5896.3335
$ws.Range("K65").Value = 22442.5
$ws.Range("L65").Value = 29481.6675
$ws.Range("M65").Value = -19322.5
$ws.Range("N65").Value = -35721.6675

$ws.Range("H95").Value = 34974.332
$ws.Range("J95").Value = 34974.332
$ws.Range("L95").Value = 34974.332
$ws.Range("N95").Value = -40466.332

$ws.Range("H102").Value = 45000
$ws.Range("J102").Value = 45000
$ws.Range("L102").Value = 45000
$ws.Range("N102").Value = -51490

$ws.Range("H132").Value = 3205.6316
$ws.Range("I132").Value = 3406.2942
$ws.Range("K132").Value = 10218.8826
$ws.Range("M132").Value = -7688.882599999999

$ws.Range("H135").Value = 521
$ws.Range("I135").Value = 252.4
$ws.Range("K135").Value = 2271.6
$ws.Range("M135").Value = 263.4000000000001

$ws.Range("H136").Value = 50000
$ws.Range("J136").Value = 50000
$ws.Range("L136").Value = 50000
$ws.Range("N136").Value = -60200

$ws.Range("H138").Value = 3469.9092
$ws.Range("I138").Value = 2194.8333
$ws.Range("K138").Value = 6584.499899999999
$ws.Range("M138").Value = -1444.499899999999

$ws.Range("H141").Value = 3359.4119
$ws.Range("J141").Value = 984
$ws.Range("L141").Value = 2952
$ws.Range("N141").Value = -13312

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 37000
$ws.Range("J24").Value = 37000
$ws.Range("L24").Value = 37000
$ws.Range("N24").Value = -37748

$ws.Range("H92").Value = 49950
$ws.Range("J92").Value = 49950
$ws.Range("L92").Value = 49950
$ws.Range("N92").Value = -54942

$ws.Range("H100").Value = 37000
$ws.Range("J100").Value = 37000
$ws.Range("L100").Value = 37000
$ws.Range("N100").Value = -39164

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 2250
$ws.Range("I5").Value = 1367
$ws.Range("J5").Value = 2912.25
$ws.Range("K5").Value = 1367
$ws.Range("L5").Value = 2912.25
$ws.Range("M5").Value = -1254
$ws.Range("N5").Value = -3138.25

$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 920
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 727.5
$ws.Range("I22").Value = 703.3333
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 703.3333
$ws.Range("L22").Value = 800
$ws.Range("M22").Value = -353.3333
$ws.Range("N22").Value = -1500

$ws.Range("H107").Value = 931.6667
$ws.Range("I107").Value = 363.33334
$ws.Range("K107").Value = 363.33334
$ws.Range("M107").Value = 1556.66666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 486
$ws.Range("I7").Value = 250.5
$ws.Range("J7").Value = 580.2
$ws.Range("K7").Value = 751.5
$ws.Range("L7").Value = 1740.6
$ws.Range("M7").Value = -639.5
$ws.Range("N7").Value = -1964.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1083.4
$ws.Range("I5").Value = 1083.4
$ws.Range("K5").Value = 1083.4
$ws.Range("M5").Value = -971.4000000000001

$ws.Range("H97").Value = 888.55554
$ws.Range("I97").Value = 567
$ws.Range("K97").Value = 567
$ws.Range("M97").Value = -71

$ws.Range("H107").Value = 838
$ws.Range("I107").Value = 300
$ws.Range("J107").Value = 1645
$ws.Range("K107").Value = 300
$ws.Range("L107").Value = 1645
$ws.Range("M107").Value = 1620
$ws.Range("N107").Value = -5485

$ws.Range("H122").Value = 4583.3335
$ws.Range("I122").Value = 2500.3333
$ws.Range("J122").Value = 6666.3335
$ws.Range("K122").Value = 7500.999899999999
$ws.Range("L122").Value = 19999.0005
$ws.Range("M122").Value = -5050.999899999999
$ws.Range("N122").Value = -24899.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1975
$ws.Range("I7").Value = 1900
$ws.Range("K7").Value = 1900
$ws.Range("M7").Value = -1788

$ws.Range("H10").Value = 1525
$ws.Range("I10").Value = 1983.3334
$ws.Range("K10").Value = 1983.3334
$ws.Range("M10").Value = -1843.3334

$ws.Range("H22").Value = 1031.6428
$ws.Range("I22").Value = 1086.625
$ws.Range("J22").Value = 958.3333
$ws.Range("K22").Value = 1086.625
$ws.Range("L22").Value = 958.3333
$ws.Range("M22").Value = -791.625
$ws.Range("N22").Value = -1548.3333

$ws.Range("H27").Value = 1031.6428
$ws.Range("I27").Value = 1086.625
$ws.Range("J27").Value = 958.3333
$ws.Range("K27").Value = 1086.625
$ws.Range("L27").Value = 958.3333
$ws.Range("M27").Value = -979.625
$ws.Range("N27").Value = -1172.3333

$ws.Range("H76").Value = 30000
$ws.Range("J76").Value = 30000
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30676

$ws.Range("H79").Value = 30000
$ws.Range("J79").Value = 30000
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32340

$ws.Range("H82").Value = 1727.4
$ws.Range("I82").Value = 996.3333
$ws.Range("J82").Value = 2040.7142
$ws.Range("K82").Value = 996.3333
$ws.Range("L82").Value = 2040.7142
$ws.Range("M82").Value = -635.3333
$ws.Range("N82").Value = -2762.7142

$ws.Range("H85").Value = 1727.4
$ws.Range("I85").Value = 996.3333
$ws.Range("J85").Value = 2040.7142
$ws.Range("K85").Value = 996.3333
$ws.Range("L85").Value = 2040.7142
$ws.Range("M85").Value = 251.6667
$ws.Range("N85").Value = -4536.7142

$ws.Range("H93").Value = 857.2
$ws.Range("I93").Value = 857.2
$ws.Range("K93").Value = 857.2
$ws.Range("M93").Value = 390.8

$ws.Range("H122").Value = 3051.4
$ws.Range("I122").Value = 3144.8572
$ws.Range("J122").Value = 2833.3333
$ws.Range("K122").Value = 9434.571599999999
$ws.Range("L122").Value = 8499.999899999999
$ws.Range("M122").Value = -6984.571599999999
$ws.Range("N122").Value = -13399.9999

$ws.Range("H126").Value = 1975
$ws.Range("I126").Value = 1900
$ws.Range("K126").Value = 5700
$ws.Range("M126").Value = -3230

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 878.7143
$ws.Range("I81").Value = 878.7143
$ws.Range("K81").Value = 1757.4286
$ws.Range("M81").Value = -696.4286

$ws.Range("H84").Value = 878.7143
$ws.Range("I84").Value = 878.7143
$ws.Range("K84").Value = 8787.143
$ws.Range("M84").Value = -3483.143

$ws.Range("H95").Value = 20000
$ws.Range("J95").Value = 20000
$ws.Range("L95").Value = 20000
$ws.Range("N95").Value = -25492

$ws.Range("H97").Value = 35373.25
$ws.Range("J97").Value = 35373.25
$ws.Range("L97").Value = 35373.25
$ws.Range("N97").Value = -37355.25

$ws.Range("H100").Value = 1800
$ws.Range("I100").Value = 1800
$ws.Range("K100").Value = 3600
$ws.Range("M100").Value = -3059

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()

$ws.Range("H113").Value = 7018.3125
$ws.Range("I113").Value = 14964
$ws.Range("K113").Value = 44892
$ws.Range("M113").Value = -42722

